$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.94975957313360326
$ws.Range("AX1").Value = 0.95124717916822277
$ws.Range("BA1").Value = 0.7338862566645068
$ws.Range("BO1").Value = 0.92259858592233934
$ws.Range("C2").Value = 0.8630633465823474
$ws.Range("D2").Value = 0.59819376128136614
$ws.Range("BP2").Value = 0.76623199824352672
$ws.Range("E3").Value = 0.90701335075009082
$ws.Range("D5").Value = 0.96424481294087583
$ws.Range("E6").Value = 0.87950733384783963
$ws.Range("H6").Value = 0.82732689200458842
$ws.Range("E7").Value = 0.89070493197407363
$ws.Range("F7").Value = 0.76918049925644949
$ws.Range("I7").Value = 0.67459534508937979
$ws.Range("G8").Value = 0.80547870237520303
$ws.Range("S8").Value = 0.98741791018325409
$ws.Range("AM8").Value = 0.66368516960571022
$ws.Range("K9").Value = 0.99499636192007634
$ws.Range("H10").Value = 0.63511219099235783
$ws.Range("I10").Value = 0.99475176926649023
$ws.Range("L10").Value = 0.64040784228808789
$ws.Range("J11").Value = 0.63233905612391683
$ws.Range("M11").Value = 0.99148301939685335
$ws.Range("K12").Value = 0.84925147818846658
$ws.Range("BI12").Value = 0.87775065159688714
$ws.Range("L13").Value = 0.97545220301880775
$ws.Range("N13").Value = 0.72365467254557025
$ws.Range("L14").Value = 0.93255538630356538
$ws.Range("M15").Value = 0.57215956754964714
$ws.Range("N15").Value = 0.97184194742359808
$ws.Range("AX15").Value = 0.92843669511774685
$ws.Range("N16").Value = 0.83416074375067717
$ws.Range("O16").Value = 0.86712428985761347
$ws.Range("R16").Value = 0.90384538516084523
$ws.Range("O17").Value = 0.92892573103983689
$ws.Range("P17").Value = 0.83775461121460082
$ws.Range("AE17").Value = 0.63958284824032297
$ws.Range("Q18").Value = 0.84233678280989199
$ws.Range("S18").Value = 0.68170651324177156
$ws.Range("T18").Value = 0.98731988729527753
$ws.Range("AW18").Value = 0.73659698203109447
$ws.Range("V20").Value = 0.97351737965439344
$ws.Range("AD20").Value = 0.77620907905036318
$ws.Range("S21").Value = 0.76810875145275159
$ws.Range("U22").Value = 0.88696829544559674
$ws.Range("X22").Value = 0.87951815068394346
$ws.Range("U23").Value = 0.76276995828424954
$ws.Range("V23").Value = 0.7316263993561346
$ws.Range("X23").Value = 0.95303860754240799
$ws.Range("Z24").Value = 0.62900802843552683
$ws.Range("W25").Value = 0.97514356832091265
$ws.Range("X25").Value = 0.89929129905176919
$ws.Range("Y26").Value = 0.87575628935775096
$ws.Range("BG26").Value = 0.92333636033785527
$ws.Range("Y27").Value = 0.82569608688177332
$ws.Range("AB27").Value = 0.8257500951028145
$ws.Range("AC27").Value = 0.64720172031940382
$ws.Range("Z28").Value = 0.70007071791873421
$ws.Range("AD28").Value = 0.63320325127126809
$ws.Range("AR28").Value = 0.71042971719664916
$ws.Range("AB29").Value = 0.93306488090657169
$ws.Range("AD29").Value = 0.85502908907268527
$ws.Range("AM29").Value = 0.95017477870606482
$ws.Range("AE30").Value = 0.83193845426726498
$ws.Range("X31").Value = 0.74875674867081343
$ws.Range("AC31").Value = 0.93831742695440756
$ws.Range("AD32").Value = 0.90647615979514673
$ws.Range("AE32").Value = 0.89394166613047576
$ws.Range("AH32").Value = 0.87429719965860586
$ws.Range("AF33").Value = 0.64734826851184657
$ws.Range("AH33").Value = 0.71151006688242835
$ws.Range("AI34").Value = 0.72659572497346137
$ws.Range("AH36").Value = 0.91738020063571435
$ws.Range("AL36").Value = 0.71925386581898598
$ws.Range("AJ37").Value = 0.92292756026168155
$ws.Range("AM37").Value = 0.80129853229409453
$ws.Range("AK38").Value = 0.97357847750208093
$ws.Range("AM38").Value = 0.91665233085125242
$ws.Range("AN38").Value = 0.51792754562914567
$ws.Range("AO39").Value = 0.91997462423543408
$ws.Range("AO40").Value = 0.72745782296963712
$ws.Range("AQ41").Value = 0.88058085090526883
$ws.Range("AG42").Value = 0.9923821193803708
$ws.Range("AN42").Value = 0.93695303957918297
$ws.Range("AO42").Value = 0.96405450333462128
$ws.Range("AQ42").Value = 0.6606593506903109
$ws.Range("AS43").Value = 0.7993310449471418
$ws.Range("AP44").Value = 0.96555719745906199
$ws.Range("AQ44").Value = 0.97709087067885814
$ws.Range("C45").Value = 0.77490952843059702
$ws.Range("AR45").Value = 0.69563017100140168
$ws.Range("AU45").Value = 0.76632375930032992
$ws.Range("AS46").Value = 0.85264067332550453
$ws.Range("AU46").Value = 0.60398925272112303
$ws.Range("AI47").Value = 0.94683360782084147
$ws.Range("AW47").Value = 0.76803461067719714
$ws.Range("AT48").Value = 0.93686570481157461
$ws.Range("AU48").Value = 0.74023160131551125
$ws.Range("AX48").Value = 0.5623790913923763
$ws.Range("BB48").Value = 0.71976540175810466
$ws.Range("D50").Value = 0.74937488491684778
$ws.Range("AW50").Value = 0.88286911681858915
$ws.Range("AL51").Value = 0.6817404436928467
$ws.Range("BJ51").Value = 0.63312449754717881
$ws.Range("BA52").Value = 0.99767406829417249
$ws.Range("AI53").Value = 0.74574671144421845
$ws.Range("AZ54").Value = 0.73872026408474745
$ws.Range("BA54").Value = 0.91189833613852
$ws.Range("BC54").Value = 0.97281959832394704
$ws.Range("BD54").Value = 0.95956060144892208
$ws.Range("BA55").Value = 0.81377556021103725
$ws.Range("BD55").Value = 0.95542653713106729
$ws.Range("BF56").Value = 0.7337735736758999
$ws.Range("BC57").Value = 0.98275765336138243
$ws.Range("BD57").Value = 0.95764717457550608
$ws.Range("BF57").Value = 0.75148258941213575
$ws.Range("BG57").Value = 0.91938371606897062
$ws.Range("G58").Value = 0.60497453162494197
$ws.Range("BF60").Value = 0.98804066731403362
$ws.Range("BG60").Value = 0.6507235837596852
$ws.Range("BJ60").Value = 0.96663425117568624
$ws.Range("BG61").Value = 0.92540015378917706
$ws.Range("BH61").Value = 0.89836058187620793
$ws.Range("E62").Value = 0.91622985281617886
$ws.Range("BK62").Value = 0.66065300255349446
$ws.Range("BI63").Value = 0.6678457716861872
$ws.Range("BM63").Value = 0.96181100265626585
$ws.Range("BK64").Value = 0.67263551280588996
$ws.Range("BM64").Value = 0.6115868137526449
$ws.Range("BN65").Value = 0.94665974932470776
$ws.Range("BO65").Value = 0.96646641529645405
$ws.Range("BL66").Value = 0.66475705213108727
$ws.Range("BN67").Value = 0.98840658821063587
$ws.Range("A68").Value = 0.85132366798797587
$ws.Range("BN68").Value = 0.780343705060625
$ws.Range("BO68").Value = 0.93796756267585779
